# Append a new row of price data to the end of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data (row 98 -> 99).
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2024-11-09 00:00:00"
$ws.Cells.Item($newRow, 2).Value = 75400
$ws.Cells.Item($newRow, 3).Value = 10530.73
$ws.Cells.Item($newRow, 4).Value = 9319.23
$ws.Cells.Item($newRow, 5).Value = 7.1792
